# GymWorkouts.xlsx update - 13/05/2018
# Appends 13 new workout-log rows (sheet rows 852-864) covering two new
# gym sessions: Friday 11-May-2018 (session 105) and Saturday 12-May-2018
# (session 106), mirroring the column layout/formatting already used by
# the existing data (A:K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows: A-value(Id), B(DateId), C(ExerciseDate serial),
# D(ExerciseMonth), E(ExerciseYear), F(ExerciseDay), G(ExerciseName),
# H(Weight), I(Sets), J(Reps), K(TrainingArea)
$newRows = @(
    @(851, 105, 43231, "May", 2018, "Friday",   "Pec Fly",         105,  4, 8,  "Chest"),
    @(852, 105, 43231, "May", 2018, "Friday",   "Hip adduction",   60,   4, 12, "Legs"),
    @(853, 105, 43231, "May", 2018, "Friday",   "Hip abduction",   60,   4, 12, "Legs"),
    @(854, 105, 43231, "May", 2018, "Friday",   "Upright Rows",    30,   4, 12, "Shoulders"),
    @(855, 105, 43231, "May", 2018, "Friday",   "Leg Extension",   101,  4, 8,  "Legs"),
    @(856, 105, 43231, "May", 2018, "Friday",   "Bicep Curl",      25,   4, 12, "Arms"),
    @(857, 106, 43232, "May", 2018, "Saturday", "Bench Press",     85,   5, 5,  "Chest"),
    @(858, 106, 43232, "May", 2018, "Saturday", "Overhead Press",  52.5, 5, 5,  "Shoulders"),
    @(859, 106, 43232, "May", 2018, "Saturday", "Barbell Row",     90,   5, 5,  "Back"),
    @(860, 106, 43232, "May", 2018, "Saturday", "Front raises",    10,   4, 8,  "Shoulders"),
    @(861, 106, 43232, "May", 2018, "Saturday", "Laterial Raises", 10,   4, 8,  "Shoulders"),
    @(862, 106, 43232, "May", 2018, "Saturday", "Rear delt flys",  10,   4, 8,  "Back"),
    @(863, 106, 43232, "May", 2018, "Saturday", "Plank",           0,    5, 30, "Core")
)

# The row immediately above (851) already carries the exact number
# formats/styles (integer / date / 2-decimal) used throughout the sheet;
# clone it down onto each new row before filling in the values so the
# new cells pick up matching styles (s="3"/"2"/"1" as appropriate).
$formatSource = $ws.Range("A851:K851")

$startRow = 852
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $destRow = $startRow + $i
    $formatSource.Copy() | Out-Null
    $ws.Range("A" + $destRow + ":K" + $destRow).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $destRow = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($destRow, 1).Value = $data[0]
    $ws.Cells.Item($destRow, 2).Value = $data[1]
    $ws.Cells.Item($destRow, 3).Value = $data[2]
    $ws.Cells.Item($destRow, 4).Value = $data[3]
    $ws.Cells.Item($destRow, 5).Value = $data[4]
    $ws.Cells.Item($destRow, 6).Value = $data[5]
    $ws.Cells.Item($destRow, 7).Value = $data[6]
    $ws.Cells.Item($destRow, 8).Value = $data[7]
    $ws.Cells.Item($destRow, 9).Value = $data[8]
    $ws.Cells.Item($destRow, 10).Value = $data[9]
    $ws.Cells.Item($destRow, 11).Value = $data[10]
}

# Match the author's final view state: window scrolled so row 850 is the
# first visible row below the frozen header, and the active cell resting
# just past the newly-added data.
$win = $excel.ActiveWindow
$win.ScrollRow = 850
$win.ScrollColumn = 1
$ws.Range("A866").Select() | Out-Null
